# Update "想去人数" (want-to-go count) values on the "展览" and "全部类型"
# sheets to match the latest generated data (output at 456a3b4).
#
# Sheet "展览"    F3: 453 -> 454
# Sheet "展览"    F4: 7   -> 9
# Sheet "全部类型" F3: 453 -> 454
# Sheet "全部类型" F4: 7   -> 9

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F3").Value = 454
    $ws.Range("F4").Value = 9
}
